$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the "Time spent" text in C8 from "3h for the moment" to "4h for the moment"
$ws.Range("C8").Value = "4h for the moment"

# Move the active selection to C8 (as recorded in the saved sheet view)
$ws.Activate()
$ws.Range("C8").Select()
